$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Controls_Synonymous")

# Update the ED_H.. labels to BS_H.. in column A (rows 2-7)
$ws2.Range("A2").Value = "BS_H26"
$ws2.Range("A3").Value = "BS_H25"
$ws2.Range("A4").Value = "BS_H37"
$ws2.Range("A5").Value = "BS_H15"
$ws2.Range("A6").Value = "BS_H20"
$ws2.Range("A7").Value = "BS_H35"

# Make Controls_Synonymous the active sheet/tab and set its selection
$ws2.Activate()
$ws2.Range("A8").Select()
